$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.169.96'
$ws.Range("E2").Value = '  -3.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.341.90'
$ws.Range("E3").Value = '  -5.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.77'
$ws.Range("E5").Value = '  -4.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.67'
$ws.Range("E6").Value = '  -2.99%  '
$ws.Range("E7").Value = '  -2.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.333.71'
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.626'
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.68'
$ws.Range("E12").Value = '  -3.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.873.52'
$ws.Range("E15").Value = '  -5.24%  '
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.339.53'
$ws.Range("E18").Value = '  -5.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.77'
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '64.093.48'
$ws.Range("E20").Value = '  -3.41%  '
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '425.65'
$ws.Range("E22").Value = '  +2.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.84'
$ws.Range("E23").Value = '  +11.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.09'
$ws.Range("E24").Value = '  -3.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.95'
$ws.Range("E25").Value = '  -2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.15'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.81'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.62'
$ws.Range("E29").Value = '  -5.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.68'
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.66'
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '595.05'
$ws.Range("E32").Value = '  -3.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.40'
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.20'
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -10.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.54'
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0754'
$ws.Range("E39").Value = '  -7.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.48'
$ws.Range("E40").Value = '  -4.91%  '
$ws.Range("E41").Value = '  -5.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.093.23'
$ws.Range("E42").Value = '  -4.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  -4.39%  '
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  -3.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("E47").Value = '  -3.94%  '
$ws.Range("E48").Value = '  -2.84%  '
$ws.Range("E49").Value = '  -4.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.21'
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.20'
$ws.Range("E51").Value = '  -5.29%  '
